# 自动更新Excel文件 - 2026-01-01 23:13:39
#
# For each product row, column E ("剩余" = days remaining) represents
# D (总天, total days of stock) minus the number of days elapsed since
# F (开始时间, start date, stored as an integer YYYYMMDD).
#
# This script advances "today" by one day (simulating the daily
# auto-update) and recomputes E for every data row accordingly. If a
# row's remaining days would drop to zero (stock exhausted), it is
# treated as restocked "today": F is reset to the new today and E is
# reset back to the full D.

function Get-OADateFromYYYYMMDD($n) {
    $s = [string]$n
    if ($s.Length -ne 8) { return $null }
    $y = [int]$s.Substring(0, 4)
    $m = [int]$s.Substring(4, 2)
    $d = [int]$s.Substring(6, 2)
    $dt = Get-Date -Year $y -Month $m -Day $d -Hour 0 -Minute 0 -Second 0
    return [math]::Floor($dt.ToOADate())
}

function Get-YYYYMMDDFromOADate($oa) {
    $epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
    $dt = $epoch.AddDays($oa)
    $y = $dt.Year
    $m = $dt.Month
    $d = $dt.Day
    $ms = [string]$m
    if ($m -lt 10) { $ms = "0$ms" }
    $ds = [string]$d
    if ($d -lt 10) { $ds = "0$ds" }
    return [int]"$y$ms$ds"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# "今天" before this update (see previous F/E values baked into the sheet)
$prevTodayOA = Get-OADateFromYYYYMMDD 20260101
# "今天" after this update - one day later, matching the commit date 2026-01-01/02
$newTodayOA = $prevTodayOA + 1

for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2

    if ($null -eq $dVal -or $null -eq $eVal -or $null -eq $fVal) { continue }

    $fOA = Get-OADateFromYYYYMMDD $fVal
    if ($null -eq $fOA) { continue }

    $dInt = [int]$dVal

    $newE = $dInt - ($newTodayOA - $fOA)

    if ($newE -le 0) {
        # Out of stock as of the new day -> restocked today.
        $newE = $dInt
        $newFOA = $newTodayOA
        $ws.Cells.Item($r, 6).Value = Get-YYYYMMDDFromOADate $newFOA
    }

    $ws.Cells.Item($r, 5).Value = $newE
}
